# Refresh the cryptos table (cols B-E, rows 2-51) with the latest scraped
# price/volume snapshot. All of these columns are plain text cells in the
# sheet (t="inlineStr"), including the "Price" column which sometimes looks
# numeric (e.g. "42.954.58"). Excel's COM layer auto-coerces a numeric-looking
# string assigned via .Value into a real number, which would flip the cell's
# stored type away from text. Prefixing the literal with a leading apostrophe
# forces Excel to keep it as text (like typing ' in the UI); ClearFormats()
# afterwards drops the transient quote-prefix/@ style Excel applies so the
# cell is left on the default (unstyled) format, matching the original file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.ClearFormats()
}

Set-TextValue 'D2' '42.954.58'
Set-TextValue 'E2' '  +0.12%  '
Set-TextValue 'D3' '2.365.97'
Set-TextValue 'E3' '  -0.42%  '
Set-TextValue 'E4' '  -0.25%  '
Set-TextValue 'D5' '319.42'
Set-TextValue 'E5' '  -3.91%  '
Set-TextValue 'D6' '107.54'
Set-TextValue 'E6' '  +5.98%  '
Set-TextValue 'E7' '  -0.64%  '
Set-TextValue 'E8' '  -0.11%  '
Set-TextValue 'D9' '0.622'
Set-TextValue 'E9' '  -1.27%  '
Set-TextValue 'D10' '41.62'
Set-TextValue 'E10' '  +1.81%  '
Set-TextValue 'D11' '0.0929'
Set-TextValue 'E11' '  +0.37%  '
Set-TextValue 'D12' '8.51'
Set-TextValue 'E12' '  +0.07%  '
Set-TextValue 'E13' '  -1.97%  '
Set-TextValue 'E14' '  +0.60%  '
Set-TextValue 'D15' '16.13'
Set-TextValue 'E15' '  -3.69%  '
Set-TextValue 'D16' '2.723.31'
Set-TextValue 'E16' '  -0.56%  '
Set-TextValue 'D17' '2.390.45'
Set-TextValue 'E17' '  -0.13%  '
Set-TextValue 'D18' '42.903.67'
Set-TextValue 'E18' '  -0.08%  '
Set-TextValue 'D19' '7.58'
Set-TextValue 'E19' '  +1.26%  '
Set-TextValue 'E20' '  -0.22%  '
Set-TextValue 'D21' '76.14'
Set-TextValue 'E21' '  -0.02%  '
Set-TextValue 'D22' '3.68'
Set-TextValue 'E22' '  -3.63%  '
Set-TextValue 'D23' '266.98'
Set-TextValue 'E23' '  -0.95%  '
Set-TextValue 'D24' '2.33'
Set-TextValue 'E24' '  -1.56%  '
Set-TextValue 'D25' '9.46'
Set-TextValue 'E25' '  -1.42%  '
Set-TextValue 'E26' '  +0.15%  '
Set-TextValue 'D27' '11.42'
Set-TextValue 'E27' '  -2.35%  '
Set-TextValue 'D28' '23.61'
Set-TextValue 'E28' '  -0.52%  '
Set-TextValue 'E29' '  +2.67%  '
Set-TextValue 'D30' '36.88'
Set-TextValue 'E30' '  +2.26%  '
Set-TextValue 'D31' '168.38'
Set-TextValue 'E31' '  -2.88%  '
Set-TextValue 'D32' '0.0899'
Set-TextValue 'E32' '  -1.73%  '
Set-TextValue 'B33' 'Filecoin'
Set-TextValue 'C33' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D33' '6.02'
Set-TextValue 'E33' '  -0.25%  '
Set-TextValue 'B34' 'WEMIXToken'
Set-TextValue 'C34' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D34' '2.91'
Set-TextValue 'E34' '  -6.51%  '
Set-TextValue 'B35' 'Kaspa'
Set-TextValue 'C35' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D35' '0.121'
Set-TextValue 'E35' '  +12.37%  '
Set-TextValue 'B36' 'Stellar'
Set-TextValue 'C36' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D36' '0.131'
Set-TextValue 'E36' '  -1.78%  '
Set-TextValue 'D37' '4.73'
Set-TextValue 'E37' '  -0.76%  '
Set-TextValue 'E38' '  +1.05%  '
Set-TextValue 'E39' '  -3.08%  '
Set-TextValue 'D40' '2.73'
Set-TextValue 'E40' '  -1.82%  '
Set-TextValue 'B41' 'BitcoinSV'
Set-TextValue 'C41' 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue 'D41' '102.10'
Set-TextValue 'E41' '  +12.17%  '
Set-TextValue 'B42' 'ARBITRUM'
Set-TextValue 'C42' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D42' '1.53'
Set-TextValue 'E42' '  +0.55%  '
Set-TextValue 'E43' '  +3.72%  '
Set-TextValue 'D44' '71.50'
Set-TextValue 'E44' '  +3.10%  '
Set-TextValue 'E45' '  -0.14%  '
Set-TextValue 'D46' '12.33'
Set-TextValue 'E46' '  +2.97%  '
Set-TextValue 'D47' '114.10'
Set-TextValue 'E47' '  -2.27%  '
Set-TextValue 'D48' '5.54'
Set-TextValue 'E48' '  +0.89%  '
Set-TextValue 'E49' '  +1.40%  '
Set-TextValue 'D50' '76.43'
Set-TextValue 'E50' '  +9.49%  '
Set-TextValue 'E51' '  +1.40%  '
